# "a little change in account"
#
# The use-case table for every scenario in this document ("ورود به حساب
# کاربری", "خروج از حساب کاربری", "ثبت نام", "ویرایش اطلاعات شخصی",
# "فراموشی رمز عبور") has a "کنشگر اصلی: " (Primary Actor) row whose
# value is "کاربر" (User). The account module now distinguishes between
# two kinds of primary actor, so every one of those rows is updated to
# "مشتری/کارمند" (Customer/Employee).
#
# Doing the replacement against the whole-document Range means it keeps
# matching (and replacing) every occurrence, i.e. behaves like "Replace
# All" across all five tables.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "کنشگر اصلی: کاربر",  # Find what
    $false,                # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                # Format
    "کنشگر اصلی: مشتری/کارمند",  # Replace with
    2                       # Replace (wdReplaceAll)
) | Out-Null

Write-Output "Done: replaced primary-actor 'کاربر' with 'مشتری/کارمند' across the document."
